# Refresh the cryptos price/volume table (Price = column D, Volume(1h) = column E).
# Note: several "Price" strings look like plain decimals (e.g. "563.14"); a bare
# assignment would make Excel auto-convert them to numbers (losing the trailing
# zeros / exact text and risking float rounding). Prefixing with a leading
# apostrophe forces Excel to keep them as literal text, matching the source data
# which stores every Price/Volume cell as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.961.20"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "2.417.41"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'563.14"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").Value = "'143.19"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").Value = "'5.20"
$ws.Range("E11").Value = "  -3.48%  "
$ws.Range("D12").Value = "'0.350"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").Value = "'25.76"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").Value = "'0.0000173"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").Value = "2.853.82"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "61.876.12"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "2.430.60"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").Value = "'323.22"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "'6.84"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("D21").Value = "'4.13"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'66.56"
$ws.Range("E23").Value = "  +2.38%  "
$ws.Range("D24").Value = "'1.74"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").Value = "'8.71"
$ws.Range("E25").Value = "  -3.31%  "
$ws.Range("D26").Value = "'561.22"
$ws.Range("E26").Value = "  -2.21%  "
$ws.Range("D27").Value = "2.537.77"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").Value = "0.0₃0939"
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("D30").Value = "'8.18"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("D31").Value = "'1.39"
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("D34").Value = "'1.50"
$ws.Range("E34").Value = "  -3.23%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'4.76"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("D37").Value = "'0.379"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").Value = "'154.06"
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("D39").Value = "'5.44"
$ws.Range("E39").Value = "  -3.32%  "
$ws.Range("D40").Value = "'18.50"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "'1.81"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("D42").Value = "'0.996"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").Value = "'2.26"
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("D44").Value = "'147.55"
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("D45").Value = "'3.63"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'0.0527"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("D47").Value = "'19.83"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("D48").Value = "'0.593"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").Value = "'0.0921"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("E51").Value = "  +0.71%  "
